$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D23 previously held an empty "/" placeholder styled s=22; it now takes on
# the schedule assignment that used to live in D25 ("Pelák Olgi / Eszenyi
# Icu"), along with D25's wrapped-text style (s=21). Copy the format from
# D25 first so the style index matches exactly, then set the value.
$ws.Range("D25").Copy()
$ws.Range("D23").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D23").Value = "Pelák Olgi`nEszenyi Icu"

# D25 keeps its existing style but gets a new schedule assignment
# ("Pelák Olgi / Tomori Marika").
$ws.Range("D25").Value = "Pelák Olgi`nTomori Marika"

# Move the active selection to D26 (was C20).
$ws.Range("D26").Select()
